$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.347
$ws.Range("B3").Value = 0.207
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 0.281
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 0.292
$ws.Range("B6").Value = 0.364
